# edit.ps1
# Applies SF4 export data to rows 11-14 (Grade 7 sections/students) as text-typed
# shared-string values, matching values produced by the original export tool.
# Excel's normal Range.Value assignment auto-converts numeric-looking strings
# ("7", "0", "50", ...) into numbers, which would store them with t="n" (or no
# type) instead of the t="s" shared-string cells the source workbook used.
# To force a literal text value while preserving the cell's existing style
# (border/format) untouched, we:
#   1. snapshot the target cell's current value into a scratch cell and copy
#      the target's formatting onto that scratch cell (so we can restore it)
#   2. flip the target cell's NumberFormat to Text ("@") so the literal
#      assignment is not reinterpreted as a number
#   3. assign the literal text value
#   4. copy the scratch cell's formatting (the target's original formatting)
#      back onto the target cell, overwriting the temporary "@" format
#   5. clear the scratch cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    $helper = $ws.Cells.Item(500, 500)
    $helper.Value = $cell.Value
    $cell.Copy() | Out-Null
    $helper.PasteSpecial(-4122) | Out-Null
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $helper.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
    $helper.Clear()
}

Set-TextValue $ws.Cells.Item(11,1) '7'
Set-TextValue $ws.Cells.Item(11,2) ' FIDELITY'
Set-TextValue $ws.Cells.Item(11,3) 'Fumino Ona Furahashi '
Set-TextValue $ws.Cells.Item(11,4) '2'
Set-TextValue $ws.Cells.Item(11,5) '0'
Set-TextValue $ws.Cells.Item(11,6) '2'
Set-TextValue $ws.Cells.Item(11,7) '1'
Set-TextValue $ws.Cells.Item(11,8) '0'
Set-TextValue $ws.Cells.Item(11,9) '1'
Set-TextValue $ws.Cells.Item(11,10) '50'
Set-TextValue $ws.Cells.Item(11,11) '0'
Set-TextValue $ws.Cells.Item(11,12) '50'
Set-TextValue $ws.Cells.Item(11,16) '0'
Set-TextValue $ws.Cells.Item(11,17) '0'
Set-TextValue $ws.Cells.Item(11,18) '0'
Set-TextValue $ws.Cells.Item(11,25) '0'
Set-TextValue $ws.Cells.Item(11,26) '0'
Set-TextValue $ws.Cells.Item(11,27) '0'
Set-TextValue $ws.Cells.Item(11,34) '2'
Set-TextValue $ws.Cells.Item(11,35) '0'
Set-TextValue $ws.Cells.Item(11,36) '2'
Set-TextValue $ws.Cells.Item(12,1) '7'
Set-TextValue $ws.Cells.Item(12,2) 'Grade 7 Oreo'
Set-TextValue $ws.Cells.Item(12,3) 'Kerby Estrella Paderogao '
Set-TextValue $ws.Cells.Item(12,4) '1'
Set-TextValue $ws.Cells.Item(12,5) '0'
Set-TextValue $ws.Cells.Item(12,6) '1'
Set-TextValue $ws.Cells.Item(12,7) '0'
Set-TextValue $ws.Cells.Item(12,8) '0'
Set-TextValue $ws.Cells.Item(12,9) '0'
Set-TextValue $ws.Cells.Item(12,10) '0'
Set-TextValue $ws.Cells.Item(12,11) '0'
Set-TextValue $ws.Cells.Item(12,12) '0'
Set-TextValue $ws.Cells.Item(12,16) '0'
Set-TextValue $ws.Cells.Item(12,17) '0'
Set-TextValue $ws.Cells.Item(12,18) '0'
Set-TextValue $ws.Cells.Item(12,25) '0'
Set-TextValue $ws.Cells.Item(12,26) '0'
Set-TextValue $ws.Cells.Item(12,27) '0'
Set-TextValue $ws.Cells.Item(12,34) '1'
Set-TextValue $ws.Cells.Item(12,35) '0'
Set-TextValue $ws.Cells.Item(12,36) '1'
Set-TextValue $ws.Cells.Item(13,1) '10'
Set-TextValue $ws.Cells.Item(13,2) ' STRAWBERRY'
Set-TextValue $ws.Cells.Item(13,3) 'None None None '
Set-TextValue $ws.Cells.Item(13,4) '0'
Set-TextValue $ws.Cells.Item(13,5) '0'
Set-TextValue $ws.Cells.Item(13,6) '0'
Set-TextValue $ws.Cells.Item(13,7) '0'
Set-TextValue $ws.Cells.Item(13,8) '0'
Set-TextValue $ws.Cells.Item(13,9) '0'
Set-TextValue $ws.Cells.Item(13,10) '0'
Set-TextValue $ws.Cells.Item(13,11) '0'
Set-TextValue $ws.Cells.Item(13,12) '0'
Set-TextValue $ws.Cells.Item(13,16) '0'
Set-TextValue $ws.Cells.Item(13,17) '0'
Set-TextValue $ws.Cells.Item(13,18) '0'
Set-TextValue $ws.Cells.Item(13,25) '0'
Set-TextValue $ws.Cells.Item(13,26) '0'
Set-TextValue $ws.Cells.Item(13,27) '0'
Set-TextValue $ws.Cells.Item(13,34) '0'
Set-TextValue $ws.Cells.Item(13,35) '0'
Set-TextValue $ws.Cells.Item(13,36) '0'
Set-TextValue $ws.Cells.Item(14,1) '10'
Set-TextValue $ws.Cells.Item(14,2) ' UNITY'
Set-TextValue $ws.Cells.Item(14,3) 'None None None '
Set-TextValue $ws.Cells.Item(14,4) '0'
Set-TextValue $ws.Cells.Item(14,5) '0'
Set-TextValue $ws.Cells.Item(14,6) '0'
Set-TextValue $ws.Cells.Item(14,7) '0'
Set-TextValue $ws.Cells.Item(14,8) '0'
Set-TextValue $ws.Cells.Item(14,9) '0'
Set-TextValue $ws.Cells.Item(14,10) '0'
Set-TextValue $ws.Cells.Item(14,11) '0'
Set-TextValue $ws.Cells.Item(14,12) '0'
Set-TextValue $ws.Cells.Item(14,16) '0'
Set-TextValue $ws.Cells.Item(14,17) '0'
Set-TextValue $ws.Cells.Item(14,18) '0'
Set-TextValue $ws.Cells.Item(14,25) '0'
Set-TextValue $ws.Cells.Item(14,26) '0'
Set-TextValue $ws.Cells.Item(14,27) '0'
Set-TextValue $ws.Cells.Item(14,34) '0'
Set-TextValue $ws.Cells.Item(14,35) '0'
Set-TextValue $ws.Cells.Item(14,36) '0'
